$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.533835
$ws.Range("H2").Value = 27.06767
$ws.Range("I2").Value = 0.2387004729612817
$ws.Range("J2").Value = 0.188326562212335
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.889399666666667
$ws.Range("N2").Value = 17.668199
$ws.Range("O2").Value = 0.8160192454225522
$ws.Range("P2").Value = 0.8160192454225521
$ws.Range("Q2").Value = 79.70616333772166
$ws.Range("R2").Value = 478.23698002633
$ws.Range("S2").Value = 0.1947841798278714
$ws.Range("T2").Value = 0.153678099189533

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.533835
$ws.Range("H3").Value = 27.06767
$ws.Range("I3").Value = 0.2387004729612817
$ws.Range("J3").Value = 0.188326562212335
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.327831666666667
$ws.Range("N3").Value = 3.983495
$ws.Range("O3").Value = 0.1839807545774479
$ws.Range("P3").Value = 0.1839807545774478
$ws.Range("Q3").Value = 17.97065468444167
$ws.Range("R3").Value = 107.82392810665
$ws.Range("S3").Value = 0.0439162931334103
$ws.Range("T3").Value = 0.03464846302280208

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.783142666666667
$ws.Range("H4").Value = 14.349428
$ws.Range("I4").Value = 0.08436178043951471
$ws.Range("J4").Value = 0.09983786727684438
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.889399666666667
$ws.Range("N4").Value = 17.668199
$ws.Range("O4").Value = 0.8160192454225522
$ws.Range("P4").Value = 0.8160192454225521
$ws.Range("Q4").Value = 28.16983882668578
$ws.Range("R4").Value = 253.528549440172
$ws.Range("S4").Value = 0.06884083641675583
$ws.Range("T4").Value = 0.08146962111984746

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.783142666666667
$ws.Range("H5").Value = 14.349428
$ws.Range("I5").Value = 0.08436178043951471
$ws.Range("J5").Value = 0.09983786727684438
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.327831666666667
$ws.Range("N5").Value = 3.983495
$ws.Range("O5").Value = 0.1839807545774479
$ws.Range("P5").Value = 0.1839807545774478
$ws.Range("Q5").Value = 6.351208298984445
$ws.Range("R5").Value = 57.16087469086001
$ws.Range("S5").Value = 0.0155209440227589
$ws.Range("T5").Value = 0.01836824615699692

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.650526666666666
$ws.Range("H6").Value = 10.95158
$ws.Range("I6").Value = 0.06438547846128644
$ws.Range("J6").Value = 0.07619693206668192
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.889399666666667
$ws.Range("N6").Value = 17.668199
$ws.Range("O6").Value = 0.8160192454225522
$ws.Range("P6").Value = 0.8160192454225521
$ws.Range("Q6").Value = 21.49941053382445
$ws.Range("R6").Value = 193.49469480442
$ws.Range("S6").Value = 0.05253978955014894
$ws.Range("T6").Value = 0.06217816300856725

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.650526666666666
$ws.Range("H7").Value = 10.95158
$ws.Range("I7").Value = 0.06438547846128644
$ws.Range("J7").Value = 0.07619693206668192
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.327831666666667
$ws.Range("N7").Value = 3.983495
$ws.Range("O7").Value = 0.1839807545774479
$ws.Range("P7").Value = 0.1839807545774478
$ws.Range("Q7").Value = 4.847284908011111
$ws.Range("R7").Value = 43.62556417210001
$ws.Range("S7").Value = 0.0118456889111375
$ws.Range("T7").Value = 0.01401876905811467

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.194588000000001
$ws.Range("H8").Value = 24.583764
$ws.Range("I8").Value = 0.1445305067870891
$ws.Range("J8").Value = 0.1710444881424727
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.889399666666667
$ws.Range("N8").Value = 17.668199
$ws.Range("O8").Value = 0.8160192454225522
$ws.Range("P8").Value = 0.8160192454225521
$ws.Range("Q8").Value = 48.26120383567068
$ws.Range("R8").Value = 434.3508345210361
$ws.Range("S8").Value = 0.1179396750889395
$ws.Range("T8").Value = 0.1395755941477072

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.194588000000001
$ws.Range("H9").Value = 24.583764
$ws.Range("I9").Value = 0.1445305067870891
$ws.Range("J9").Value = 0.1710444881424727
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.327831666666667
$ws.Range("N9").Value = 3.983495
$ws.Range("O9").Value = 0.1839807545774479
$ws.Range("P9").Value = 0.1839807545774478
$ws.Range("Q9").Value = 10.88103344168667
$ws.Range("R9").Value = 97.92930097518001
$ws.Range("S9").Value = 0.0265908316981496
$ws.Range("T9").Value = 0.03146889399476545

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 13.70308866666667
$ws.Range("H10").Value = 41.10926600000001
$ws.Range("I10").Value = 0.2416856527188127
$ws.Range("J10").Value = 0.2860226514085782
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.889399666666667
$ws.Range("N10").Value = 17.668199
$ws.Range("O10").Value = 0.8160192454225522
$ws.Range("P10").Value = 0.8160192454225521
$ws.Range("Q10").Value = 80.70296582577045
$ws.Range("R10").Value = 726.3266924319341
$ws.Range("S10").Value = 0.1972201439610626
$ws.Range("T10").Value = 0.2333999881761857

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.70308866666667
$ws.Range("H11").Value = 41.10926600000001
$ws.Range("I11").Value = 0.2416856527188127
$ws.Range("J11").Value = 0.2860226514085782
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.327831666666667
$ws.Range("N11").Value = 3.983495
$ws.Range("O11").Value = 0.1839807545774479
$ws.Range("P11").Value = 0.1839807545774478
$ws.Range("Q11").Value = 18.19539506274111
$ws.Range("R11").Value = 163.75855556467
$ws.Range("S11").Value = 0.04446550875775018
$ws.Range("T11").Value = 0.05262266323239255

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 12.8328005
$ws.Range("H12").Value = 25.665601
$ws.Range("I12").Value = 0.2263361086320154
$ws.Range("J12").Value = 0.1785714988930879
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 5.889399666666667
$ws.Range("N12").Value = 17.668199
$ws.Range("O12").Value = 0.8160192454225522
$ws.Range("P12").Value = 0.8160192454225521
$ws.Range("Q12").Value = 75.57749098709984
$ws.Range("R12").Value = 453.4649459225991
$ws.Range("S12").Value = 0.1846946205777741
$ws.Range("T12").Value = 0.1457177797807117

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 12.8328005
$ws.Range("H13").Value = 25.665601
$ws.Range("I13").Value = 0.2263361086320154
$ws.Range("J13").Value = 0.1785714988930879
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.327831666666667
$ws.Range("N13").Value = 3.983495
$ws.Range("O13").Value = 0.1839807545774479
$ws.Range("P13").Value = 0.1839807545774478
$ws.Range("Q13").Value = 17.03979887591584
$ws.Range("R13").Value = 102.238793255495
$ws.Range("S13").Value = 0.04164148805424141
$ws.Range("T13").Value = 0.0328537191123762

